# Commiting new test cases
# Change the Runmode column (C) for TestCase_B2..TestCase_B6 (rows 3-7) on the
# "Test Cases" sheet from "N" to "Y", and update the active selection to
# reflect the range that was edited (C2:C7, active cell C2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("C3:C7").Value = "Y"

$ws.Range("C2:C7").Select()
